$d = $word.ActiveDocument
$p = $d.Paragraphs(1)

# Append " V" as its own bold run at the end of the first paragraph,
# without swallowing the paragraph mark (so it lands in this paragraph,
# not the start of the next one).
$r = $p.Range
$r.MoveEnd(1, -1)
$r.Collapse(0)
$r.InsertAfter(" V")
$r.Font.Bold = $true
$r.Font.BoldBi = $true

# Append "4" as another separate bold run right after it.
$p2 = $d.Paragraphs(1)
$r2 = $p2.Range
$r2.MoveEnd(1, -1)
$r2.Collapse(0)
$r2.InsertAfter("4")
$r2.Font.Bold = $true
$r2.Font.BoldBi = $true
